# Applies the commit "fin rapport histoire des sciences":
#   1. Adds a due-date reminder paragraph (+ blank line) at the top of
#      the document.
#   2. Appends a new "Consignes Maumy" section (intro line + a bulleted
#      list) right after the "... difficultés envisagées." paragraph,
#      before the trailing blank paragraphs / section break.

$d = $word.ActiveDocument

# --- Insertion 1: new paragraphs at the very start of the document ---
$startRange = $d.Range(0, 0)
$topXml = '<w:p><w:r><w:t>Remise du dossier le lundi 04/01/21 à 12h maxi.</w:t></w:r></w:p><w:p/>'
$startRange.InsertXML($topXml)

# --- Insertion 2: big block of new paragraphs, right after the
#     "... difficultés envisagées." paragraph, before the trailing
#     empty paragraphs / section break ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.Contains("envisag")) {
        $target = $para
    }
}
if ($target -eq $null) {
    throw "Could not locate target paragraph (difficultes envisagees)."
}
$insertPos = $target.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$bigXml = '<w:p/><w:p><w:r><w:t xml:space="preserve">Consignes </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Maumy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> :</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Au S3, élaborer une activité la préparer. Au S4 faire cette activité en classe et analyser le déroulement de l’activité.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Élaborer : donner les prérequis (ceux du collège et lycée), l’inclure dans la progression, définir les objectifs disciplinaires et </w:t></w:r><w:r><w:t>en termes de</w:t></w:r><w:r><w:t xml:space="preserve"> compétences. </w:t></w:r><w:r><w:t>Cette activité doit être évaluée.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Introduction : présenter notre situation de professeur stagiaire, </w:t></w:r><w:r><w:t xml:space="preserve">du lycée, de l’emploi du temps, du nombre d’élèves, de classes, </w:t></w:r><w:r><w:t>présentation éventuelle des CSP</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">1 Présentation en détail de la séance, la replacer dans le cadre de la séquence, décortiquer la séance (il faut donner le </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>déroulé</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> de la séance, minuté), les prérequis, les objectifs programmes et compétences, les difficultés éventuelles des élèves</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Évidemment joindre le document annexe.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Justifier les choix pédagogiques (travail en îlot) </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Présentation orale de cette séance.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Au S4 : rendre compte du déroulement réel de la séance</w:t></w:r></w:p>'
$insertRange.InsertXML($bigXml)

Write-Host "Edit applied."
